$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.966.33'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '2.043.29'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''245.39'
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("D6").Value = '''0.658'
$ws.Range("D7").Value = '''58.20'
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '''0.377'
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").Value = '''0.0768'
$ws.Range("E10").Value = '  -0.86%  '
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("D12").Value = '''15.49'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").Value = '''0.881'
$ws.Range("E13").Value = '  +9.38%  '
$ws.Range("D14").Value = '2.341.00'
$ws.Range("D15").Value = '''5.66'
$ws.Range("E15").Value = '  +2.40%  '
$ws.Range("D16").Value = '2.030.43'
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("D17").Value = '''18.19'
$ws.Range("E17").Value = '  +10.00%  '
$ws.Range("D18").Value = '36.957.83'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").Value = '''73.77'
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("D20").Value = '0.0₃0889'
$ws.Range("E20").Value = '  -1.25%  '
$ws.Range("D21").Value = '''5.39'
$ws.Range("E21").Value = '  +1.35%  '
$ws.Range("D22").Value = '''235.84'
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +3.50%  '
$ws.Range("D25").Value = '''9.59'
$ws.Range("E25").Value = '  +5.56%  '
$ws.Range("D26").Value = '''169.60'
$ws.Range("E26").Value = '  +1.36%  '
$ws.Range("E27").Value = '  -2.51%  '
$ws.Range("D28").Value = '''19.93'
$ws.Range("E28").Value = '  +1.29%  '
$ws.Range("D29").Value = '''5.38'
$ws.Range("E29").Value = '  +15.27%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("D31").Value = '''1.11'
$ws.Range("E31").Value = '  -1.41%  '
$ws.Range("D32").Value = '''4.70'
$ws.Range("E32").Value = '  +5.96%  '
$ws.Range("D33").Value = '''0.0612'
$ws.Range("E33").Value = '  +0.42%  '
$ws.Range("D34").Value = '''0.999'
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '''0.0872'
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("D36").Value = '''1.85'
$ws.Range("E36").Value = '  +6.63%  '
$ws.Range("D37").Value = '''2.24'
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("E39").Value = '  -2.19%  '
$ws.Range("D40").Value = '''5.09'
$ws.Range("E40").Value = '  +3.85%  '
$ws.Range("D41").Value = '''0.0983'
$ws.Range("E41").Value = '  -7.03%  '
$ws.Range("D42").Value = '''0.0222'
$ws.Range("E42").Value = '  +1.13%  '
$ws.Range("D43").Value = '''1.14'
$ws.Range("E43").Value = '  +1.75%  '
$ws.Range("D44").Value = '''96.95'
$ws.Range("E44").Value = '  +1.79%  '
$ws.Range("D45").Value = '''16.86'
$ws.Range("E45").Value = '  -2.24%  '
$ws.Range("D46").Value = '1.294.18'
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("D47").Value = '''3.81'
$ws.Range("E47").Value = '  +10.43%  '
$ws.Range("D48").Value = '''2.33'
$ws.Range("E48").Value = '  -3.31%  '
$ws.Range("D49").Value = '''2.85'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").Value = '''6.73'
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("D51").Value = '2.226.29'
$ws.Range("E51").Value = '  -0.10%  '
